# Duplicate the "function_sheet" worksheet into a new worksheet named
# "first_test", placed immediately after the original (so tab order is
# function_sheet, first_test), while keeping "function_sheet" the active
# / selected sheet -- mirrors Excel's "Move or Copy… > Create a copy"
# placed to the right of the source sheet.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(1)

# Create the new sheet right after the source sheet.
$new = $wb.Worksheets.Add($null, $src)
$new.Name = "first_test"

# Copy the data + formatting of the used range onto the new sheet.
$src.Range("A1:F6").Copy($new.Range("A1:F6"))

# Mirror the source sheet's column widths on the new sheet.
for ($c = 1; $c -le 6; $c++) {
    $new.Columns.Item($c).ColumnWidth = $src.Columns.Item($c).ColumnWidth
}

# Keep the original sheet as the active / selected tab.
$src.Activate()
